$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 561
$ws.Range("J18").Value = 340
$ws.Range("L18").Value = 340
$ws.Range("N18").Value = -908
$ws.Range("H64").Value = 3828.1667
$ws.Range("I64").Value = 3942
$ws.Range("J64").Value = 3600.5
$ws.Range("K64").Value = 3942
$ws.Range("L64").Value = 3600.5
$ws.Range("M64").Value = -3694
$ws.Range("N64").Value = -4096.5
$ws.Range("H67").Value = 3828.1667
$ws.Range("I67").Value = 3942
$ws.Range("J67").Value = 3600.5
$ws.Range("K67").Value = 3942
$ws.Range("L67").Value = 3600.5
$ws.Range("M67").Value = -3084
$ws.Range("N67").Value = -5316.5
$ws.Range("H116").Value = 7633.3335
$ws.Range("I116").Value = 10700
$ws.Range("K116").Value = 10700
$ws.Range("M116").Value = -7258
$ws.Range("H123").Value = 39419.8
$ws.Range("J123").Value = 39419.8
$ws.Range("L123").Value = 39419.8
$ws.Range("N123").Value = -49219.8
$ws.Range("H126").Value = 35685
$ws.Range("J126").Value = 35685
$ws.Range("L126").Value = 35685
$ws.Range("N126").Value = -45565
$ws.Range("H127").Value = 1841.75
$ws.Range("I127").Value = 591.4
$ws.Range("J127").Value = 2043.4193
$ws.Range("K127").Value = 1774.2
$ws.Range("L127").Value = 6130.257900000001
$ws.Range("M127").Value = 3185.8
$ws.Range("N127").Value = -16050.2579
$ws.Range("H137").Value = 1426.069
$ws.Range("I137").Value = 1258.762
$ws.Range("J137").Value = 1865.25
$ws.Range("K137").Value = 3776.286
$ws.Range("L137").Value = 5595.75
$ws.Range("M137").Value = -1226.286
$ws.Range("N137").Value = -10695.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4876.279
$ws.Range("I32").Value = 3545.6667
$ws.Range("J32").Value = 11662.4
$ws.Range("K32").Value = 3545.6667
$ws.Range("L32").Value = 11662.4
$ws.Range("M32").Value = -3258.6667
$ws.Range("N32").Value = -12236.4
$ws.Range("H61").Value = 6143.773
$ws.Range("I61").Value = 6341.095
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 6341.095
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -6129.095
$ws.Range("N61").Value = -2424
$ws.Range("H63").Value = 111113496
$ws.Range("I63").Value = 111113496
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 111113496
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -111112810
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 111113496
$ws.Range("I66").Value = 111113496
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 555567480
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -555564048
$ws.Range("N66").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 49900
$ws.Range("J109").Value = 49900
$ws.Range("L109").Value = 49900
$ws.Range("N109").Value = -52674
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
$ws.Range("H136").Value = 6143.773
$ws.Range("I136").Value = 6341.095
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 19023.285
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -16473.285
$ws.Range("N136").Value = -11100
$ws.Range("H141").Value = 39429
$ws.Range("J141").Value = 39429
$ws.Range("L141").Value = 39429
$ws.Range("N141").Value = -49789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 28471
$ws.Range("I26").Value = 28471
$ws.Range("K26").Value = 28471
$ws.Range("M26").Value = -28179
$ws.Range("H59").Value = 49780
$ws.Range("J59").Value = 49780
$ws.Range("L59").Value = 49780
$ws.Range("N59").Value = -51474
$ws.Range("H60").Value = 39780
$ws.Range("J60").Value = 39780
$ws.Range("L60").Value = 39780
$ws.Range("N60").Value = -40978
$ws.Range("H74").Value = 41306.715
$ws.Range("J74").Value = 47758.5
$ws.Range("L74").Value = 47758.5
$ws.Range("N74").Value = -49630.5
$ws.Range("H77").Value = 41306.715
$ws.Range("J77").Value = 47758.5
$ws.Range("L77").Value = 143275.5
$ws.Range("N77").Value = -152635.5
$ws.Range("H81").Value = 35900
$ws.Range("J81").Value = 35900
$ws.Range("L81").Value = 35900
$ws.Range("N81").Value = -38022
$ws.Range("H84").Value = 35900
$ws.Range("J84").Value = 35900
$ws.Range("L84").Value = 107700
$ws.Range("N84").Value = -118308
$ws.Range("H122").Value = 44499
$ws.Range("J122").Value = 44499
$ws.Range("L122").Value = 44499
$ws.Range("N122").Value = -54299
$ws.Range("H126").Value = 47765.77
$ws.Range("J126").Value = 47765.77
$ws.Range("L126").Value = 47765.77
$ws.Range("N126").Value = -57645.77
$ws.Range("H139").Value = 50317.25
$ws.Range("J139").Value = 60186.668
$ws.Range("L139").Value = 60186.668
$ws.Range("N139").Value = -70466.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 6833.3335
$ws.Range("I41").Value = 2750
$ws.Range("K41").Value = 2750
$ws.Range("M41").Value = -2322
$ws.Range("H50").Value = 21295
$ws.Range("J50").Value = 21295
$ws.Range("L50").Value = 21295
$ws.Range("N50").Value = -22545
$ws.Range("H51").Value = 24894
$ws.Range("J51").Value = 24894
$ws.Range("L51").Value = 24894
$ws.Range("N51").Value = -26366
$ws.Range("H59").Value = 21142
$ws.Range("J59").Value = 21142
$ws.Range("L59").Value = 21142
$ws.Range("N59").Value = -23432
$ws.Range("H60").Value = 18883.857
$ws.Range("I60").Value = 6546.5
$ws.Range("J60").Value = 23818.8
$ws.Range("K60").Value = 6546.5
$ws.Range("L60").Value = 23818.8
$ws.Range("M60").Value = -6035.5
$ws.Range("N60").Value = -24840.8
$ws.Range("H61").Value = 24894
$ws.Range("J61").Value = 24894
$ws.Range("L61").Value = 24894
$ws.Range("N61").Value = -25590
$ws.Range("H62").Value = 4756.6665
$ws.Range("I62").Value = 5153.077
$ws.Range("J62").Value = 4112.5
$ws.Range("K62").Value = 5153.077
$ws.Range("L62").Value = 4112.5
$ws.Range("M62").Value = -4529.077
$ws.Range("N62").Value = -5360.5
$ws.Range("H65").Value = 4756.6665
$ws.Range("I65").Value = 5153.077
$ws.Range("J65").Value = 4112.5
$ws.Range("K65").Value = 25765.385
$ws.Range("L65").Value = 20562.5
$ws.Range("M65").Value = -22645.385
$ws.Range("N65").Value = -26802.5
$ws.Range("H68").Value = 27225
$ws.Range("J68").Value = 27225
$ws.Range("L68").Value = 27225
$ws.Range("N68").Value = -28723
$ws.Range("H71").Value = 27225
$ws.Range("J71").Value = 27225
$ws.Range("L71").Value = 81675
$ws.Range("N71").Value = -89163
$ws.Range("H74").Value = 21932.834
$ws.Range("J74").Value = 21932.834
$ws.Range("L74").Value = 21932.834
$ws.Range("N74").Value = -23680.834
$ws.Range("H77").Value = 21932.834
$ws.Range("J77").Value = 21932.834
$ws.Range("L77").Value = 65798.50199999999
$ws.Range("N77").Value = -74534.50199999999
$ws.Range("H132").Value = 2874.7666
$ws.Range("I132").Value = 2988.682
$ws.Range("J132").Value = 2561.5
$ws.Range("K132").Value = 8966.045999999998
$ws.Range("L132").Value = 7684.5
$ws.Range("M132").Value = -6436.045999999998
$ws.Range("N132").Value = -12744.5
$ws.Range("H135").Value = 36140.625
$ws.Range("J135").Value = 36140.625
$ws.Range("L135").Value = 36140.625
$ws.Range("N135").Value = -46280.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 34195
$ws.Range("J133").Value = 34195
$ws.Range("L133").Value = 34195
$ws.Range("N133").Value = -44315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1737505.8
$ws.Range("I22").Value = 4630101
$ws.Range("J22").Value = 1948.575
$ws.Range("K22").Value = 4630101
$ws.Range("L22").Value = 1948.575
$ws.Range("M22").Value = -4629806
$ws.Range("N22").Value = -2538.575
$ws.Range("H27").Value = 1737505.8
$ws.Range("I27").Value = 4630101
$ws.Range("J27").Value = 1948.575
$ws.Range("K27").Value = 4630101
$ws.Range("L27").Value = 1948.575
$ws.Range("M27").Value = -4629994
$ws.Range("N27").Value = -2162.575
$ws.Range("H40").Value = 76925180
$ws.Range("I40").Value = 111112880
$ws.Range("J40").Value = 2873.75
$ws.Range("K40").Value = 111112880
$ws.Range("L40").Value = 2873.75
$ws.Range("M40").Value = -111112744
$ws.Range("N40").Value = -3145.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820
